$wb = $excel.ActiveWorkbook

# 1. Update "Metadata" last-updated timestamp (A2)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "18 Nov 2025, 09:20 AM"

# 2. "Top Losers" sheet: rows 17-76 shift up by two positions
#    (TMCV / RAJRILTD drop off the top, NATIONALUM / AHLUCONT are
#    appended at the bottom) - refresh Stock/Latest/Weekly/Monthly.
$losers = $wb.Worksheets.Item("Top Losers")

$loserData = @(
    @(17, "VIDHIING", -1.5303, -2.9498, 7.5751),
    @(18, "WINDMACHIN", -1.5278, -0.7867, -1.0068),
    @(19, "SEAMECLTD", -1.5143, -0.2882, -3.7659),
    @(20, "CARYSIL", -1.4302, -1.5492, 9.038399999999999),
    @(21, "PAYTM", -1.418, 1.1392, 0.8287),
    @(22, "MPHASIS", -1.4022, -4.3685, -4.1022),
    @(23, "RAYMONDLSL", -1.3934, -1.2005, -6.0306),
    @(24, "THEJO", -1.3801, -1.5259, -7.4692),
    @(25, "BHAGCHEM", -1.2951, 0.5179, 0.3028),
    @(26, "SIGNPOST", -1.281, -3.7431, -4.6405),
    @(27, "GODAVARIB", -1.2746, -2.2749, -5.613),
    @(28, "HEMIPROP", -1.2721, -1.1491, 2.9857),
    @(29, "JARO", -1.1718, 11.6249, 11.1464),
    @(30, "EKC", -1.1097, -7.5512, -10.4947),
    @(31, "SESHAPAPER", -1.1065, -2.5316, -2.0164),
    @(32, "CREDITACC", -1.1062, -1.8804, -5.8022),
    @(33, "PRIMESECU", -1.1005, -1.3993, 6.8401),
    @(34, "BECTORFOOD", -1.0751, 7.6623, -1.0825),
    @(35, "ALICON", -1.0731, 0.6139, -4.7055),
    @(36, "ASHIANA", -1.0577, -1.8613, 2.7822),
    @(37, "BFUTILITIE", -1.0262, 0.0519, -8.3367),
    @(38, "AKUMS", -1.0216, 2.7575, -5.3717),
    @(39, "HINDCOPPER", -1.019, -0.0447, -1.5136),
    @(40, "IPCALAB", -1.0104, -2.0005, 14.0655),
    @(41, "ANURAS", -1.0038, 0.1677, 0.0745),
    @(42, "INDOSTAR", -0.9935, 0.7112000000000001, -0.4317),
    @(43, "PRECAM", -0.9836, 1.3477, 1.9396),
    @(44, "GIPCL", -0.9666, 1.8949, -5.2945),
    @(45, "INDIQUBE", -0.9131, -1.6138, 1.7534),
    @(46, "STERTOOLS", -0.9071, 0.704, -3.4859),
    @(47, "JASH", -0.9054, -2.8962, -1.926),
    @(48, "NRBBEARING", -0.9012, -0.418, 4.4193),
    @(49, "CHEMPLASTS", -0.8988, -2.55, -15.6826),
    @(50, "SWELECTES", -0.8955, -4.8529, -7.1358),
    @(51, "GUFICBIO", -0.8937, -2.4199, -4.579),
    @(52, "STEELCAS", -0.8875999999999999, -0.9873, -3.6816),
    @(53, "VPRPL", -0.8871, -3.4242, -4.8269),
    @(54, "BAJAJINDEF", -0.8857, -1.3906, -6.4891),
    @(55, "INNOVACAP", -0.8807, -3.3741, -12.6519),
    @(56, "SOLARA", -0.8522999999999999, -0.5149, -4.9683),
    @(57, "NAVNETEDUL", -0.8453000000000001, 1.114, -5.6607),
    @(58, "PANAMAPET", -0.8386, -0.2239, 8.1172),
    @(59, "METROPOLIS", -0.8223, -0.2436, -0.2233),
    @(60, "HINDZINC", -0.806, -1.4475, 0.7345),
    @(61, "KIOCL", -0.7887999999999999, -2.8665, -13.8596),
    @(62, "CAMPUS", -0.7618, 0.3117, -2.2512),
    @(63, "EUROPRATIK", -0.7596000000000001, 7.2156, 19.1959),
    @(64, "GREENPLY", -0.7568, -0.8902, -3.4206),
    @(65, "BSOFT", -0.7292, -0.1801, 4.4133),
    @(66, "SANGHIIND", -0.729, -1.036, -1.5687),
    @(67, "GVPIL", -0.7244, 4.24, 21.7706),
    @(68, "AMRUTANJAN", -0.7232, -1.7851, -2.3615),
    @(69, "GPTINFRA", -0.7166, -1.5059, -1.3676),
    @(70, "FUSION", -0.7141999999999999, 0.0351, -6.0801),
    @(71, "DATAMATICS", -0.7045, 1.8029, -4.2012),
    @(72, "WEL", -0.6849, 3.5714, 11.7088),
    @(73, "HUHTAMAKI", -0.6764, -1.9649, -5.128),
    @(74, "DHANUKA", -0.6753, -0.1579, -7.9537),
    @(75, "NATIONALUM", -0.675, -1.3596, 10.6176),
    @(76, "AHLUCONT", -0.6724, 10.4315, 2.5722)
)

foreach ($row in $loserData) {
    $r = $row[0]
    $losers.Cells.Item($r, 2).Value = $row[1]
    $losers.Cells.Item($r, 3).Value = $row[2]
    $losers.Cells.Item($r, 4).Value = $row[3]
    $losers.Cells.Item($r, 5).Value = $row[4]
}

# 3. "Industry Analysis" sheet: F7 (textiles - processing) refresh
$industry = $wb.Worksheets.Item("Industry Analysis")
$industry.Cells.Item(7, 6).Value = 7.5375

